$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(99, 8).Value = 226.85715
$ws.Cells.Item(99, 10).Value = 400
$ws.Cells.Item(99, 12).Value = 1200
$ws.Cells.Item(99, 14).Value = -4196
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 32574.727
$ws.Cells.Item(32, 9).Value = 5387.4165
$ws.Cells.Item(32, 10).Value = 219002
$ws.Cells.Item(32, 11).Value = 5387.4165
$ws.Cells.Item(32, 12).Value = 219002
$ws.Cells.Item(32, 13).Value = -5100.4165
$ws.Cells.Item(32, 14).Value = -219576
$ws.Cells.Item(74, 8).Value = 2482.3333
$ws.Cells.Item(74, 9).Value = 827.5714
$ws.Cells.Item(74, 11).Value = 827.5714
$ws.Cells.Item(74, 13).Value = 46.42859999999996
$ws.Cells.Item(77, 8).Value = 2482.3333
$ws.Cells.Item(77, 9).Value = 827.5714
$ws.Cells.Item(77, 11).Value = 4137.857
$ws.Cells.Item(77, 13).Value = 230.143
$ws.Cells.Item(139, 8).Value = 33571.668
$ws.Cells.Item(139, 10).Value = 33571.668
$ws.Cells.Item(139, 12).Value = 33571.668
$ws.Cells.Item(139, 14).Value = -43851.668
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(99, 8).Value = 1483.9615
$ws.Cells.Item(99, 9).Value = 1215.6666
$ws.Cells.Item(99, 10).Value = 1713.9286
$ws.Cells.Item(99, 11).Value = 1215.6666
$ws.Cells.Item(99, 12).Value = 1713.9286
$ws.Cells.Item(99, 13).Value = 282.3334
$ws.Cells.Item(99, 14).Value = -4709.9286
$ws.Cells.Item(105, 8).Value = 89334.35000000001
$ws.Cells.Item(105, 9).Value = 65074.938
$ws.Cells.Item(105, 10).Value = 144784.42
$ws.Cells.Item(105, 11).Value = 65074.938
$ws.Cells.Item(105, 12).Value = 144784.42
$ws.Cells.Item(105, 13).Value = -63327.938
$ws.Cells.Item(105, 14).Value = -148278.42
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 40479.22
$ws.Cells.Item(31, 9).Value = 1415.9166
$ws.Cells.Item(31, 10).Value = 56643.344
$ws.Cells.Item(31, 11).Value = 1415.9166
$ws.Cells.Item(31, 12).Value = 56643.344
$ws.Cells.Item(31, 13).Value = -1120.9166
$ws.Cells.Item(31, 14).Value = -57233.344
$ws.Cells.Item(34, 8).Value = 40479.22
$ws.Cells.Item(34, 9).Value = 1415.9166
$ws.Cells.Item(34, 10).Value = 56643.344
$ws.Cells.Item(34, 11).Value = 1415.9166
$ws.Cells.Item(34, 12).Value = 56643.344
$ws.Cells.Item(34, 13).Value = -1213.9166
$ws.Cells.Item(34, 14).Value = -57047.344
$ws.Cells.Item(58, 8).Value = 3923.8845
$ws.Cells.Item(58, 9).Value = 937.587
$ws.Cells.Item(58, 10).Value = 26818.834
$ws.Cells.Item(58, 11).Value = 937.587
$ws.Cells.Item(58, 12).Value = 26818.834
$ws.Cells.Item(58, 13).Value = -734.587
$ws.Cells.Item(58, 14).Value = -27224.834
$ws.Cells.Item(105, 8).Value = 2527.2727
$ws.Cells.Item(105, 9).Value = 2698.625
$ws.Cells.Item(105, 10).Value = 2070.3333
$ws.Cells.Item(105, 11).Value = 2698.625
$ws.Cells.Item(105, 12).Value = 2070.3333
$ws.Cells.Item(105, 13).Value = -951.625
$ws.Cells.Item(105, 14).Value = -5564.3333
$ws.Cells.Item(136, 8).Value = 3923.8845
$ws.Cells.Item(136, 9).Value = 937.587
$ws.Cells.Item(136, 10).Value = 26818.834
$ws.Cells.Item(136, 11).Value = 2812.761
$ws.Cells.Item(136, 12).Value = 80456.50199999999
$ws.Cells.Item(136, 13).Value = -262.761
$ws.Cells.Item(136, 14).Value = -85556.50199999999
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(22, 8).Value = 8189.1113
$ws.Cells.Item(22, 9).Value = 500
$ws.Cells.Item(22, 10).Value = 9150.25
$ws.Cells.Item(22, 11).Value = 1500
$ws.Cells.Item(22, 12).Value = 27450.75
$ws.Cells.Item(22, 13).Value = -1331
$ws.Cells.Item(22, 14).Value = -27788.75
$ws.Cells.Item(27, 8).Value = 8189.1113
$ws.Cells.Item(27, 9).Value = 500
$ws.Cells.Item(27, 10).Value = 9150.25
$ws.Cells.Item(27, 11).Value = 1500
$ws.Cells.Item(27, 12).Value = 27450.75
$ws.Cells.Item(27, 13).Value = -1398
$ws.Cells.Item(27, 14).Value = -27654.75
$ws.Cells.Item(34, 8).Value = 720.3
$ws.Cells.Item(34, 10).Value = 977.5714
$ws.Cells.Item(34, 12).Value = 2932.7142
$ws.Cells.Item(34, 14).Value = -3100.7142
$ws.Cells.Item(40, 8).Value = 421.42856
$ws.Cells.Item(40, 9).Value = 313.625
$ws.Cells.Item(40, 10).Value = 487.76923
$ws.Cells.Item(40, 11).Value = 1254.5
$ws.Cells.Item(40, 12).Value = 1951.07692
$ws.Cells.Item(40, 13).Value = -1185.5
$ws.Cells.Item(40, 14).Value = -2089.07692
$ws.Cells.Item(46, 8).Value = 37192.855
$ws.Cells.Item(46, 9).Value = 1600
$ws.Cells.Item(46, 10).Value = 46900
$ws.Cells.Item(46, 11).Value = 4800
$ws.Cells.Item(46, 12).Value = 140700
$ws.Cells.Item(46, 13).Value = -4709
$ws.Cells.Item(46, 14).Value = -140882
$ws.Cells.Item(131, 8).Value = 1092.45
$ws.Cells.Item(131, 9).Value = 1030
$ws.Cells.Item(131, 10).Value = 1093.0808
$ws.Cells.Item(131, 11).Value = 3090
$ws.Cells.Item(131, 12).Value = 3279.2424
$ws.Cells.Item(131, 13).Value = 1950
$ws.Cells.Item(131, 14).Value = -13359.2424
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(69, 8).Value = 0
$ws.Cells.Item(69, 10).Value = 0
$ws.Cells.Item(69, 12).Value = 0
$ws.Cells.Item(69, 14).ClearContents()
$ws.Cells.Item(72, 8).Value = 0
$ws.Cells.Item(72, 10).Value = 0
$ws.Cells.Item(72, 12).Value = 0
$ws.Cells.Item(72, 14).ClearContents()
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(16, 8).Value = 32118.594
$ws.Cells.Item(16, 9).Value = 37710.965
$ws.Cells.Item(16, 10).Value = 1919.8
$ws.Cells.Item(16, 11).Value = 37710.965
$ws.Cells.Item(16, 12).Value = 1919.8
$ws.Cells.Item(16, 13).Value = -37540.965
$ws.Cells.Item(16, 14).Value = -2259.8
$ws.Cells.Item(46, 8).Value = 596068.3
$ws.Cells.Item(46, 9).Value = 419.6
$ws.Cells.Item(46, 10).Value = 844255.25
$ws.Cells.Item(46, 11).Value = 419.6
$ws.Cells.Item(46, 12).Value = 844255.25
$ws.Cells.Item(46, 13).Value = -231.6
$ws.Cells.Item(46, 14).Value = -844631.25
$ws.Cells.Item(100, 8).Value = 2338.3635
$ws.Cells.Item(100, 9).Value = 1974.75
$ws.Cells.Item(100, 10).Value = 2546.1428
$ws.Cells.Item(100, 11).Value = 1974.75
$ws.Cells.Item(100, 12).Value = 2546.1428
$ws.Cells.Item(100, 13).Value = -1433.75
$ws.Cells.Item(100, 14).Value = -3628.1428
$ws.Cells.Item(132, 8).Value = 2591.2766
$ws.Cells.Item(132, 9).Value = 2712.0264
$ws.Cells.Item(132, 10).Value = 2081.4443
$ws.Cells.Item(132, 11).Value = 8136.0792
$ws.Cells.Item(132, 12).Value = 6244.3329
$ws.Cells.Item(132, 13).Value = -5606.0792
$ws.Cells.Item(132, 14).Value = -11304.3329
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(100, 8).Value = 48514.43
$ws.Cells.Item(100, 9).Value = 84368.664
$ws.Cells.Item(100, 10).Value = 708.7778
$ws.Cells.Item(100, 11).Value = 168737.328
$ws.Cells.Item(100, 12).Value = 1417.5556
$ws.Cells.Item(100, 13).Value = -168196.328
$ws.Cells.Item(100, 14).Value = -2499.5556
$ws.Cells.Item(132, 8).Value = 2106.9106
$ws.Cells.Item(132, 9).Value = 1972.3914
$ws.Cells.Item(132, 11).Value = 5917.174199999999
$ws.Cells.Item(132, 13).Value = -3387.174199999999
$ws.Cells.Item(136, 8).Value = 948.9677
$ws.Cells.Item(136, 9).Value = 612.5714
$ws.Cells.Item(136, 10).Value = 1655.4
$ws.Cells.Item(136, 11).Value = 1837.7142
$ws.Cells.Item(136, 12).Value = 4966.200000000001
$ws.Cells.Item(136, 13).Value = 712.2857999999999
$ws.Cells.Item(136, 14).Value = -10066.2
